# The paragraph "Er moet een systeem in komen waarmee een admin permissies
# kan geven aan andere docenten, vergelijkbaar met hoe het in ELO gebeurt."
# (together with the _GoBack bookmark that follows it) is merged into the
# next paragraph, which is the one carrying the section's sectPr. This is
# equivalent to deleting the paragraph mark at the end of that paragraph.

$d = $word.ActiveDocument

$target = "Er moet een systeem in komen waarmee een admin permissies kan geven aan andere docenten, vergelijkbaar met hoe het in ELO gebeurt."

$r = $d.Content
$found = $r.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph text to merge."
}

# Collapse the found range to its end (wdCollapseEnd = 0), then extend it
# by one character (wdCharacter = 1) so it spans the paragraph mark that
# ends the paragraph, and delete that mark. This merges the paragraph with
# the following one (the empty, sectPr-carrying paragraph), exactly like
# pressing Delete at the end of the line in Word.
$r.Collapse(0) | Out-Null
$r.MoveEnd(1, 1) | Out-Null
$r.Delete()
